# Commit: "correções dos calculos ok"
#  - TESTES sheet: the "Todos elementos negativos" test now passes ("ok"),
#    and a new test row is added for "Primeiro elemento com sinal negativo"
#    (also "ok"), right below it.
#  - PROBLEMAS sheet: cursor left parked on B3 after review.

$wb = $excel.ActiveWorkbook

# --- TESTES sheet ------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("TESTES")

# Row 3 ("Todos elementos negativos") flips from "não ok" to "ok".
$ws1.Range("B3").Value = "ok"

# Insert a new row 5 for the negative-first-element test case, pushing the
# existing rows 5 and 6 down to 6 and 7.
$ws1.Rows.Item(5).Insert()
$ws1.Range("A5").Value = "Primeiro elemento com sinal negativo"
$ws1.Range("B5").Value = "ok"

# Column B needs to be wide enough to fit "ok"/"não ok" (matches the
# bestFit width Excel would compute for this column; 5.5 is the closest
# achievable ColumnWidth to the target stored width of 6.25).
$ws1.Columns.Item(2).ColumnWidth = 5.5

$ws1.Activate()
$ws1.Range("B7").Select()

# --- PROBLEMAS sheet ---------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PROBLEMAS")
$ws2.Range("B3").Select()

# Leave TESTES as the active/selected sheet & cell, matching the final state.
$ws1.Activate()
$ws1.Range("B7").Select()
